$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.04529833333333
$ws.Range("H2").Value = 45.135895
$ws.Range("I2").Value = 0.5098953433306311
$ws.Range("J2").Value = 0.5098953433306311
$ws.Range("M2").Value = 10.82167433333333
$ws.Range("N2").Value = 32.465023
$ws.Range("O2").Value = 0.09133543757015983
$ws.Range("P2").Value = 0.09133543757015983
$ws.Range("Q2").Value = 162.8153188111761
$ws.Range("R2").Value = 1465.337869300585
$ws.Range("S2").Value = 0.04657151429809007
$ws.Range("T2").Value = 0.04657151429809007

# Row 3
$ws.Range("G3").Value = 15.04529833333333
$ws.Range("H3").Value = 45.135895
$ws.Range("I3").Value = 0.5098953433306311
$ws.Range("J3").Value = 0.5098953433306311
$ws.Range("M3").Value = 36.14140700000001
$ws.Range("O3").Value = 0.3050351656377608
$ws.Range("P3").Value = 0.3050351656377608
$ws.Range("Q3").Value = 543.7582505014218
$ws.Range("R3").Value = 4893.824254512796
$ws.Range("S3").Value = 0.155536010510782
$ws.Range("T3").Value = 0.155536010510782

# Row 4
$ws.Range("G4").Value = 15.04529833333333
$ws.Range("H4").Value = 45.135895
$ws.Range("I4").Value = 0.5098953433306311
$ws.Range("J4").Value = 0.5098953433306311
$ws.Range("M4").Value = 26.40107466666666
$ws.Range("N4").Value = 79.20322399999999
$ws.Range("O4").Value = 0.2228263051286729
$ws.Range("P4").Value = 0.2228263051286729
$ws.Range("Q4").Value = 397.2120446806089
$ws.Range("R4").Value = 3574.90840212548
$ws.Range("S4").Value = 0.1136180953566806
$ws.Range("T4").Value = 0.1136180953566806

# Row 5
$ws.Range("G5").Value = 15.04529833333333
$ws.Range("H5").Value = 45.135895
$ws.Range("I5").Value = 0.5098953433306311
$ws.Range("J5").Value = 0.5098953433306311
$ws.Range("M5").Value = 45.11859966666666
$ws.Range("N5").Value = 135.355799
$ws.Range("O5").Value = 0.3808030916634065
$ws.Range("P5").Value = 0.3808030916634065
$ws.Range("Q5").Value = 678.8227923672339
$ws.Range("R5").Value = 6109.405131305105
$ws.Range("S5").Value = 0.1941697231650784
$ws.Range("T5").Value = 0.1941697231650785

# Row 6
$ws.Range("I6").Value = 0.2313044792629727
$ws.Range("J6").Value = 0.2313044792629727
$ws.Range("M6").Value = 10.82167433333333
$ws.Range("N6").Value = 32.465023
$ws.Range("O6").Value = 0.09133543757015983
$ws.Range("P6").Value = 0.09133543757015983
$ws.Range("Q6").Value = 73.858122115138
$ws.Range("R6").Value = 664.723099036242
$ws.Range("S6").Value = 0.02112629582542157
$ws.Range("T6").Value = 0.02112629582542157

# Row 7
$ws.Range("I7").Value = 0.2313044792629727
$ws.Range("J7").Value = 0.2313044792629727
$ws.Range("M7").Value = 36.14140700000001
$ws.Range("O7").Value = 0.3050351656377608
$ws.Range("P7").Value = 0.3050351656377608
$ws.Range("S7").Value = 0.07055600014473687
$ws.Range("T7").Value = 0.07055600014473687

# Row 8
$ws.Range("I8").Value = 0.2313044792629727
$ws.Range("J8").Value = 0.2313044792629727
$ws.Range("M8").Value = 26.40107466666666
$ws.Range("N8").Value = 79.20322399999999
$ws.Range("O8").Value = 0.2228263051286729
$ws.Range("P8").Value = 0.2228263051286729
$ws.Range("Q8").Value = 180.187809819344
$ws.Range("R8").Value = 1621.690288374096
$ws.Range("S8").Value = 0.05154072247387993
$ws.Range("T8").Value = 0.05154072247387994

# Row 9
$ws.Range("I9").Value = 0.2313044792629727
$ws.Range("J9").Value = 0.2313044792629727
$ws.Range("M9").Value = 45.11859966666666
$ws.Range("N9").Value = 135.355799
$ws.Range("O9").Value = 0.3808030916634065
$ws.Range("P9").Value = 0.3808030916634065
$ws.Range("Q9").Value = 307.935254859794
$ws.Range("R9").Value = 2771.417293738146
$ws.Range("S9").Value = 0.08808146081893428
$ws.Range("T9").Value = 0.0880814608189343

# Row 10
$ws.Range("G10").Value = 6.326195000000001
$ws.Range("H10").Value = 18.978585
$ws.Range("I10").Value = 0.2143990301843924
$ws.Range("J10").Value = 0.2143990301843924
$ws.Range("M10").Value = 10.82167433333333
$ws.Range("N10").Value = 32.465023
$ws.Range("O10").Value = 0.09133543757015983
$ws.Range("P10").Value = 0.09133543757015983
$ws.Range("Q10").Value = 68.46002205916169
$ws.Range("R10").Value = 616.1401985324551
$ws.Range("S10").Value = 0.01958222923650938
$ws.Range("T10").Value = 0.01958222923650938

# Row 11
$ws.Range("G11").Value = 6.326195000000001
$ws.Range("H11").Value = 18.978585
$ws.Range("I11").Value = 0.2143990301843924
$ws.Range("J11").Value = 0.2143990301843924
$ws.Range("M11").Value = 36.14140700000001
$ws.Range("O11").Value = 0.3050351656377608
$ws.Range("P11").Value = 0.3050351656377608
$ws.Range("Q11").Value = 228.6375882563651
$ws.Range("R11").Value = 2057.738294307286
$ws.Range("S11").Value = 0.06539924368487142
$ws.Range("T11").Value = 0.06539924368487141

# Row 12
$ws.Range("G12").Value = 6.326195000000001
$ws.Range("H12").Value = 18.978585
$ws.Range("I12").Value = 0.2143990301843924
$ws.Range("J12").Value = 0.2143990301843924
$ws.Range("M12").Value = 26.40107466666666
$ws.Range("N12").Value = 79.20322399999999
$ws.Range("O12").Value = 0.2228263051286729
$ws.Range("P12").Value = 0.2228263051286729
$ws.Range("Q12").Value = 167.0183465508933
$ws.Range("R12").Value = 1503.16511895804
$ws.Range("S12").Value = 0.04777374371915897
$ws.Range("T12").Value = 0.04777374371915897

# Row 13
$ws.Range("G13").Value = 6.326195000000001
$ws.Range("H13").Value = 18.978585
$ws.Range("I13").Value = 0.2143990301843924
$ws.Range("J13").Value = 0.2143990301843924
$ws.Range("M13").Value = 45.11859966666666
$ws.Range("N13").Value = 135.355799
$ws.Range("O13").Value = 0.3808030916634065
$ws.Range("P13").Value = 0.3808030916634065
$ws.Range("Q13").Value = 285.4290596182684
$ws.Range("R13").Value = 2568.861536564415
$ws.Range("S13").Value = 0.08164381354385264
$ws.Range("T13").Value = 0.08164381354385264

# Row 14
$ws.Range("G14").Value = 1.310128666666667
$ws.Range("H14").Value = 3.930386
$ws.Range("I14").Value = 0.04440114722200381
$ws.Range("J14").Value = 0.04440114722200381
$ws.Range("M14").Value = 10.82167433333333
$ws.Range("N14").Value = 32.465023
$ws.Range("O14").Value = 0.09133543757015983
$ws.Range("P14").Value = 0.09133543757015983
$ws.Range("Q14").Value = 14.17778576543089
$ws.Range("R14").Value = 127.600071888878
$ws.Range("S14").Value = 0.004055398210138805
$ws.Range("T14").Value = 0.004055398210138804

# Row 15
$ws.Range("G15").Value = 1.310128666666667
$ws.Range("H15").Value = 3.930386
$ws.Range("I15").Value = 0.04440114722200381
$ws.Range("J15").Value = 0.04440114722200381
$ws.Range("M15").Value = 36.14140700000001
$ws.Range("O15").Value = 0.3050351656377608
$ws.Range("P15").Value = 0.3050351656377608
$ws.Range("Q15").Value = 47.34989336436735
$ws.Range("R15").Value = 426.1490402793061
$ws.Range("S15").Value = 0.01354391129737053
$ws.Range("T15").Value = 0.01354391129737053

# Row 16
$ws.Range("G16").Value = 1.310128666666667
$ws.Range("H16").Value = 3.930386
$ws.Range("I16").Value = 0.04440114722200381
$ws.Range("J16").Value = 0.04440114722200381
$ws.Range("M16").Value = 26.40107466666666
$ws.Range("N16").Value = 79.20322399999999
$ws.Range("O16").Value = 0.2228263051286729
$ws.Range("P16").Value = 0.2228263051286729
$ws.Range("Q16").Value = 34.58880475160711
$ws.Range("R16").Value = 311.2992427644639
$ws.Range("S16").Value = 0.009893743578953347
$ws.Range("T16").Value = 0.009893743578953347

# Row 17
$ws.Range("G17").Value = 1.310128666666667
$ws.Range("H17").Value = 3.930386
$ws.Range("I17").Value = 0.04440114722200381
$ws.Range("J17").Value = 0.04440114722200381
$ws.Range("M17").Value = 45.11859966666666
$ws.Range("N17").Value = 135.355799
$ws.Range("O17").Value = 0.3808030916634065
$ws.Range("P17").Value = 0.3808030916634065
$ws.Range("Q17").Value = 59.1111708231571
$ws.Range("R17").Value = 532.0005374084139
$ws.Range("S17").Value = 0.01690809413554113
$ws.Range("T17").Value = 0.01690809413554113
